$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 127-128; this shifts the existing rows 127..241
# down to 129..243 (matching the target dimension A1:R243) and carries
# the column-D date style (s="2") into the new rows automatically.
$ws.Rows("127:128").Insert()

# Row 127 - new record (Primera)
$ws.Range("A127").Value = 1
$ws.Range("B127").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C127").Value = "Arica y Parinacota"
$ws.Range("D127").Value = 44586
$ws.Range("E127").Value = 15
$ws.Range("F127").Value = 100112043
$ws.Range("G127").Value = "Pepino ensalada"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 120
$ws.Range("K127").Value = 7000
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = 7250
$ws.Range("N127").Value = "$/caja 70 unidades"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 104
$ws.Range("Q127").Value = 70
$ws.Range("R127").Value = "Hortaliza"

# Row 128 - new record (Segunda)
$ws.Range("A128").Value = 1
$ws.Range("B128").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C128").Value = "Arica y Parinacota"
$ws.Range("D128").Value = 44586
$ws.Range("E128").Value = 15
$ws.Range("F128").Value = 100112043
$ws.Range("G128").Value = "Pepino ensalada"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Segunda"
$ws.Range("J128").Value = 120
$ws.Range("K128").Value = 6000
$ws.Range("L128").Value = 6500
$ws.Range("M128").Value = 6250
$ws.Range("N128").Value = "$/caja 100 unidades"
$ws.Range("O128").Value = "Región de Arica y Parinacota"
$ws.Range("P128").Value = 62
$ws.Range("Q128").Value = 100
$ws.Range("R128").Value = "Hortaliza"
